$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "PAN" header (column H) to "Pan"
$ws.Range("H1").Value = "Pan"

# Add two new trailing columns: "DP" (N) and "Client Id" (O)
$ws.Range("N1").Value = "DP"
$ws.Range("O1").Value = "Client Id"

# Match the header formatting used by the neighbouring "Update Only" column (M)
$ws.Range("N1").Style = $ws.Range("M1").Style
$ws.Range("O1").Style = $ws.Range("M1").Style

# Populate the new "DP" values
$ws.Range("N2").Value = 123456
$ws.Range("N3").Value = 234567
$ws.Range("N4").Value = 345678
$ws.Range("N5").Value = 456789
$ws.Range("N6").Value = 567900
$ws.Range("N7").Value = 679011

# Populate the new "Client Id" values
$ws.Range("O2").Value = 1
$ws.Range("O3").Value = 2
$ws.Range("O4").Value = 3
$ws.Range("O5").Value = 4
$ws.Range("O6").Value = 5
$ws.Range("O7").Value = 6

# Move the selection to reflect where the user ended up after entering the data
$ws.Range("O8").Select()
